# Append " (Changed main)" -- as three distinct runs -- right after the
# existing "This is a Microsoft word document." run in the first paragraph.
$d = $word.ActiveDocument

$rng = $d.Content
$rng.Find.Execute("This is a Microsoft word document.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# $rng now spans exactly the found text. Re-insert that same text plus the
# three new trailing runs as an OOXML fragment so each piece of text lands
# in its own <w:r>, instead of Word's usual "extend the existing run"
# behaviour when calling InsertAfter with matching formatting.
$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>This is a Microsoft word document.</w:t></w:r><w:r><w:t xml:space="preserve"> (</w:t></w:r><w:r><w:t>Changed main</w:t></w:r><w:r><w:t>)</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$rng.InsertXML($xml)
